$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row: insert new "Jan_2026" column before "Dec_2025",
# shifting existing month columns and dropping the oldest month (Oct_2025).
$ws.Range("D1").Value = "Jan_2026"
$ws.Range("E1").Value = "Dec_2025"
$ws.Range("F1").Value = "Nov_2025"
$ws.Range("G1").Value = "MoM"
$ws.Range("H1").Value = "QoQ"

# Full refreshed data set (rows 2-22), re-ranked by the new Jan_2026 holding value.
$data = @(
    @("INE040A01034", "HDFC Bank Limited", "quant Aggressive Hybrid Fund", 8.303181, 0.486552, 0, 7.816629000000001, 8.303181),
    @("INE018A01030", "Larsen & Toubro Limited", "quant Aggressive Hybrid Fund", 7.967097, 7.795835, 7.650868, 0.1712619999999996, 0.3162289999999999),
    @("INE814H01029", "Adani Power Limited", "quant Aggressive Hybrid Fund", 7.270346, 7.235738, 7.349683, 0.03460800000000042, -0.07933699999999977),
    @("INE917I01010", "Bajaj Auto Limited", "quant Aggressive Hybrid Fund", 7.01815, 6.437652, 5.181581, 0.5804980000000004, 1.836569),
    @("INE090A01021", "ICICI Bank Limited", "quant Aggressive Hybrid Fund", 6.16377, 0, 0, 6.16377, 6.16377),
    @("INE364U01010", "Adani Green Energy Limited", "quant Aggressive Hybrid Fund", 5.924793, 6.64915, 3.708483, -0.7243569999999995, 2.21631),
    @("INE795G01014", "HDFC Life Insurance Co Ltd", "quant Aggressive Hybrid Fund", 5.775625, 5.58255, 5.603391, 0.1930749999999994, 0.1722339999999996),
    @("INE556F16BR2", "SIDBI CD 10-Nov-2026", "quant Aggressive Hybrid Fund", 4.834687, 4.551729, 0, 0.2829579999999998, 4.834687),
    @("INE423A01024", "Adani Enterprises Limited", "quant Aggressive Hybrid Fund", 4.516326, 4.717528, 3.448005, -0.2012019999999994, 1.068321),
    @("INE775A01035", "Samvardhana Motherson International Ltd", "quant Aggressive Hybrid Fund", 3.673839, 3.67762, 0, -0.003781000000000034, 3.673839),
    @("INE406A01037", "Aurobindo Pharma Limited", "quant Aggressive Hybrid Fund", 3.444555, 3.13663, 3.202917, 0.307925, 0.241638),
    @("INE758E01017", "Jio Financial Services Limited", "quant Aggressive Hybrid Fund", 3.106554, 7.429747, 7.595557, -4.323193, -4.489003),
    @("INE127D01025", "HDFC Asset Management Company Ltd", "quant Aggressive Hybrid Fund", 2.975697, 0, 0, 2.975697, 2.975697),
    @("INE0Z8Z23013", "Capital Infra Trust InvIT", "quant Aggressive Hybrid Fund", 1.895889, 1.810264, 1.825513, 0.08562499999999984, 0.07037599999999999),
    @("INE721A01047", "Shriram Finance Limited", "quant Aggressive Hybrid Fund", 1.450518, 0, 0, 1.450518, 1.450518),
    @("INE245A01021", "Tata Power Company Limited", "quant Aggressive Hybrid Fund", 0, 3.840354, 5.288741, -3.840354, -5.288741),
    @("INE423A20016", "Adani Enterprises Limited Rights", "quant Aggressive Hybrid Fund", 0, 0, 0.080276, 0, -0.080276),
    @("INE216A01030", "Britannia Industries Limited", "quant Aggressive Hybrid Fund", 0, 6.408499, 6.117228, -6.408499, -6.117228),
    @("INE669C01036", "Tech Mahindra Limited", "quant Aggressive Hybrid Fund", 0, 0.335577, 0, -0.335577, 0),
    @("INE062A01020", "State Bank of India", "quant Aggressive Hybrid Fund", 0, 0, 5.281295, 0, -5.281295),
    @("INE002A01018", "Reliance Industries Limited", "quant Aggressive Hybrid Fund", 0, 5.255267, 8.443042, -5.255267, -8.443042)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
    $row++
}
